$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Excused Folio Ids" column header in J1 (next empty column after the
# existing "Investment Domicile *" header in I1), matching the commit's intent of
# adding a folio-ids column to the Portfolio Investments import template.
$ws.Range("J1").Value = "Excused Folio Ids"

# The newly typed header cell picks up a plain Times New Roman 10pt font (distinct
# from the bold/bordered Arial used by the other header cells).
$ws.Range("J1").Font.Name = "Times New Roman"
$ws.Range("J1").Font.Size = 10
$ws.Range("J1").Font.Bold = $false

# Leave the selection where the editor ended up after adding the column.
$ws.Range("K5").Select()
